# Penalty Reward System update
# Shifts the forecast dates forward by one week and updates the
# MyForecast values accordingly, then refreshes the dependent
# summary statistics on the "Summary" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Forecast Comparison"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# Force the Week_Start_Date column to be treated as plain text so Excel
# does not reinterpret the date-like strings as date serial numbers.
$ws1.Range("B2:B17").NumberFormat = "@"

$weekDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

$myForecast = @(161, 140, 134, 141, 162, 192, 222, 156, 152, 152, 221, 196, 171, 160, 165, 179)

for ($i = 0; $i -lt $weekDates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).Value = $weekDates[$i]
    $ws1.Cells.Item($row, 4).Value = $myForecast[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Summary"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

# Keep these as plain text, matching the original inline-string storage.
$ws2.Range("B2").NumberFormat = "@"
$ws2.Range("B4:B15").NumberFormat = "@"

$ws2.Range("B2").Value = "2022-12-25 to 2025-01-05"
$ws2.Range("B4").Value = "293"
$ws2.Range("B5").Value = "146"
$ws2.Range("B6").Value = "146"
$ws2.Range("B8").Value = "15075 units"
$ws2.Range("B9").Value = "2704"
$ws2.Range("B10").Value = "1308"
$ws2.Range("B11").Value = "576"
$ws2.Range("B12").Value = "222"
$ws2.Range("B14").Value = "134"
